$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update. Values are written with a leading apostrophe
# marker (classic "force text" Excel convention) so number-looking
# strings like "97.53" are not silently coerced into floating point
# values - the source data are text cells (t="inlineStr") that must
# keep their exact literal formatting (trailing zeros, dot-grouped
# thousands, percent strings with padding spaces, etc).
# The Style reset afterwards clears the incidental "Text" number
# format that the apostrophe marker leaves behind, restoring the
# cell to the workbook default style so only the value changes.

$ws.Range('D2').Value = '''42.916.11'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''  -0.38%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''2.300.94'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  -0.11%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = '''  -0.04%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''305.64'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  +1.65%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''97.53'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '''  -0.10%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').Value = '''  -1.77%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('E8').Value = '''  -0.05%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('E9').Value = '''  -1.97%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = '''35.75'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  +0.05%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = '''0.0792'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '''  +0.28%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = '''18.18'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''  +1.01%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = '''0.119'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '''  +1.01%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('E14').Value = '''  -1.33%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''2.659.13'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  -0.12%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = '''2.300.30'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''  -1.26%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('E17').Value = '''  +0.01%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''42.840.03'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  -0.35%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = '''12.74'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '''  -3.60%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''0.0₃0904'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  -0.29%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = '''6.05'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '''  -0.94%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''67.78'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  -0.95%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = '''236.75'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '''  -0.56%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = '''2.16'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '''  -1.98%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''2.46'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  +1.61%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('E26').Value = '''  +0.06%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('E27').Value = '''  -0.03%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = '''25.44'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '''  +2.83%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = '''166.90'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '''  -0.60%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = '''2.05'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '''  +0.89%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = '''9.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '''  -1.09%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''33.06'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  +1.00%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('E33').Value = '''  +0.05%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''4.81'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  -0.11%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('E35').Value = '''  -2.62%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = '''17.22'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '''  -4.88%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('E37').Value = '''  -0.88%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = '''0.0692'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '''  +0.59%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('E39').Value = '''  -0.93%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('E40').Value = '''  -1.51%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('E41').Value = '''  -1.25%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = '''2.74'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''  -0.49%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = '''2.007.37'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '''  -0.06%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('E44').Value = '''  -2.00%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('B45').Value = '''FraxShare'
$ws.Range('B45').Style = "Normal"
$ws.Range('C45').Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C45').Style = "Normal"
$ws.Range('D45').Value = '''10.01'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''  -1.52%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('B46').Value = '''ApeXProtocol'
$ws.Range('B46').Style = "Normal"
$ws.Range('C46').Value = '''https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('C46').Style = "Normal"
$ws.Range('D46').Value = '''2.12'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '''  -2.13%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = '''17.96'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '''  +4.10%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = '''  -1.30%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('E49').Value = '''  +3.95%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = '''54.27'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '''  -0.07%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = '''2.526.34'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '''  +0.04%  '
$ws.Range('E51').Style = "Normal"
